# Auto-generated: re-applies the scheduled market-data refresh for the Leve profit sheets.
# For each changed cell we set the new currentAveragePrice* / LevePrice* / LeveProfit*
# value; cells that the refresh removed entirely are cleared so the row shape matches.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 539.3333
$ws.Range("I4").Value = 577.8333
$ws.Range("K4").Value = 577.8333
$ws.Range("M4").Value = -463.8333
$ws.Range("H6").Value = 91375.37
$ws.Range("I6").Value = 100492.9
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 301478.7
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -301366.7
$ws.Range("N6").Value = -824
$ws.Range("H12").Value = 10890807
$ws.Range("I12").Value = 14520744
$ws.Range("K12").Value = 14520744
$ws.Range("M12").Value = -14520574
$ws.Range("H38").Value = 615.2
$ws.Range("J38").Value = 2498
$ws.Range("L38").Value = 7494
$ws.Range("N38").Value = -8238
$ws.Range("H58").Value = 2237.25
$ws.Range("J58").Value = 4250
$ws.Range("L58").Value = 12750
$ws.Range("N58").Value = -13050
$ws.Range("H76").Value = 2965.8333
$ws.Range("I76").Value = 2149
$ws.Range("K76").Value = 2149
$ws.Range("M76").Value = -1834
$ws.Range("H79").Value = 2965.8333
$ws.Range("I79").Value = 2149
$ws.Range("K79").Value = 2149
$ws.Range("M79").Value = -1057
$ws.Range("H86").Value = 70177900
$ws.Range("I86").Value = 71431390
$ws.Range("J86").Value = 66668120
$ws.Range("K86").Value = 71431390
$ws.Range("L86").Value = 66668120
$ws.Range("M86").Value = -71430267
$ws.Range("N86").Value = -66670366
$ws.Range("H89").Value = 70177900
$ws.Range("I89").Value = 71431390
$ws.Range("J89").Value = 66668120
$ws.Range("K89").Value = 357156950
$ws.Range("L89").Value = 333340600
$ws.Range("M89").Value = -357151334
$ws.Range("N89").Value = -333351832
$ws.Range("H96").Value = 1569.7273
$ws.Range("I96").Value = 956.3333
$ws.Range("J96").Value = 2305.8
$ws.Range("K96").Value = 2868.9999
$ws.Range("L96").Value = 6917.400000000001
$ws.Range("M96").Value = -1495.9999
$ws.Range("N96").Value = -9663.400000000001
$ws.Range("H97").Value = 2797.5
$ws.Range("J97").Value = 2797.5
$ws.Range("L97").Value = 8392.5
$ws.Range("N97").Value = -9384.5
$ws.Range("H101").Value = 3857.1428
$ws.Range("I101").Value = 2333.3333
$ws.Range("K101").Value = 6999.999899999999
$ws.Range("M101").Value = -5377.999899999999
$ws.Range("H135").Value = 624.5862
$ws.Range("I135").Value = 399.7619
$ws.Range("K135").Value = 3597.8571
$ws.Range("M135").Value = -1062.8571
$ws.Range("H137").Value = 997.3043
$ws.Range("I137").Value = 937.45
$ws.Range("K137").Value = 2812.35
$ws.Range("M137").Value = -262.3500000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 94769.63
$ws.Range("I45").Value = 114940.664
$ws.Range("K45").Value = 114940.664
$ws.Range("M45").Value = -114563.664
$ws.Range("H97").Value = 7599
$ws.Range("I97").Value = 8922.333000000001
$ws.Range("K97").Value = 8922.333000000001
$ws.Range("M97").Value = -8426.333000000001
$ws.Range("H102").Value = 2486.8
$ws.Range("I102").Value = 2486.8
$ws.Range("K102").Value = 2486.8
$ws.Range("M102").Value = -864.8000000000002
$ws.Range("H132").Value = 2132.907
$ws.Range("I132").Value = 1771.5405
$ws.Range("K132").Value = 5314.6215
$ws.Range("M132").Value = -2784.6215
$ws.Range("H135").Value = 87900.336
$ws.Range("J135").Value = 87900.336
$ws.Range("L135").Value = 87900.336
$ws.Range("N135").Value = -98040.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6781.4585
$ws.Range("I99").Value = 7687.75
$ws.Range("K99").Value = 7687.75
$ws.Range("M99").Value = -6189.75
$ws.Range("H105").Value = 5143.3213
$ws.Range("J105").Value = 3406.75
$ws.Range("L105").Value = 3406.75
$ws.Range("N105").Value = -6900.75
$ws.Range("H107").Value = 13288.1
$ws.Range("I107").Value = 21761.818
$ws.Range("J107").Value = 2931.3333
$ws.Range("K107").Value = 21761.818
$ws.Range("L107").Value = 2931.3333
$ws.Range("M107").Value = -19841.818
$ws.Range("N107").Value = -6771.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2579.85
$ws.Range("I31").Value = 2712.1765
$ws.Range("J31").Value = 2482.0435
$ws.Range("K31").Value = 2712.1765
$ws.Range("L31").Value = 2482.0435
$ws.Range("M31").Value = -2417.1765
$ws.Range("N31").Value = -3072.0435
$ws.Range("H34").Value = 2579.85
$ws.Range("I34").Value = 2712.1765
$ws.Range("J34").Value = 2482.0435
$ws.Range("K34").Value = 2712.1765
$ws.Range("L34").Value = 2482.0435
$ws.Range("M34").Value = -2510.1765
$ws.Range("N34").Value = -2886.0435
$ws.Range("H43").Value = 17079.8
$ws.Range("J43").Value = 17079.8
$ws.Range("L43").Value = 17079.8
$ws.Range("N43").Value = -17447.8
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("H101").Value = 17079.8
$ws.Range("J101").Value = 17079.8
$ws.Range("L101").Value = 17079.8
$ws.Range("N101").Value = -23569.8
$ws.Range("H105").Value = 2089.5334
$ws.Range("I105").Value = 1361.9166
$ws.Range("K105").Value = 1361.9166
$ws.Range("M105").Value = 385.0834
$ws.Range("H107").Value = 1057.32
$ws.Range("I107").Value = 1234.2941
$ws.Range("K107").Value = 1234.2941
$ws.Range("M107").Value = 685.7058999999999
$ws.Range("H132").Value = 10449851
$ws.Range("I132").Value = 43261.457
$ws.Range("J132").Value = 41669616
$ws.Range("K132").Value = 129784.371
$ws.Range("L132").Value = 125008848
$ws.Range("M132").Value = -127254.371
$ws.Range("N132").Value = -125013908
$ws.Range("H134").Value = 1666.421
$ws.Range("I134").Value = 1303.3125
$ws.Range("K134").Value = 3909.9375
$ws.Range("M134").Value = -1374.9375
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 645.5
$ws.Range("I11").Value = 630.5714
$ws.Range("J11").Value = 750
$ws.Range("K11").Value = 1891.7142
$ws.Range("L11").Value = 2250
$ws.Range("M11").Value = -1751.7142
$ws.Range("N11").Value = -2530
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9000
$ws.Range("J70").Value = 9000
$ws.Range("L70").Value = 9000
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 9000
$ws.Range("J73").Value = 9000
$ws.Range("L73").Value = 9000
$ws.Range("N73").Value = -10872
$ws.Range("H126").Value = 4808.273
$ws.Range("I126").Value = 4599.8
$ws.Range("J126").Value = 4982
$ws.Range("K126").Value = 13799.4
$ws.Range("L126").Value = 14946
$ws.Range("M126").Value = -11329.4
$ws.Range("N126").Value = -19886
$ws.Range("H132").Value = 4631397.5
$ws.Range("I132").Value = 1821.2142
$ws.Range("J132").Value = 9617095
$ws.Range("K132").Value = 5463.642599999999
$ws.Range("L132").Value = 28851285
$ws.Range("M132").Value = -2933.642599999999
$ws.Range("N132").Value = -28856345

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9801.883
$ws.Range("J7").Value = 12095.546
$ws.Range("L7").Value = 12095.546
$ws.Range("N7").Value = -12319.546
$ws.Range("H9").Value = 4149.8
$ws.Range("I9").Value = 5333.3335
$ws.Range("J9").Value = 2374.5
$ws.Range("K9").Value = 5333.3335
$ws.Range("L9").Value = 2374.5
$ws.Range("M9").Value = -5109.3335
$ws.Range("N9").Value = -2822.5
$ws.Range("H40").Value = 4703.5557
$ws.Range("I40").Value = 5749.75
$ws.Range("K40").Value = 5749.75
$ws.Range("M40").Value = -5613.75
$ws.Range("H46").Value = 8495.9
$ws.Range("I46").Value = 9931.25
$ws.Range("K46").Value = 9931.25
$ws.Range("M46").Value = -9743.25
$ws.Range("H82").Value = 2680.0908
$ws.Range("I82").Value = 1957.6
$ws.Range("K82").Value = 1957.6
$ws.Range("M82").Value = -1596.6
$ws.Range("H85").Value = 2680.0908
$ws.Range("I85").Value = 1957.6
$ws.Range("K85").Value = 1957.6
$ws.Range("M85").Value = -709.5999999999999
$ws.Range("H126").Value = 9801.883
$ws.Range("J126").Value = 12095.546
$ws.Range("L126").Value = 36286.638
$ws.Range("N126").Value = -41226.638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 258.9
$ws.Range("I113").Value = 261.45834
$ws.Range("K113").Value = 784.3750200000001
$ws.Range("M113").Value = 1385.62498
$ws.Range("H132").Value = 1926.3549
$ws.Range("I132").Value = 1941.0834
$ws.Range("J132").Value = 1875.8572
$ws.Range("K132").Value = 5823.2502
$ws.Range("L132").Value = 5627.571599999999
$ws.Range("M132").Value = -3293.2502
$ws.Range("N132").Value = -10687.5716

